$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.4723746666666667
$ws.Cells.Item(2, 8).Value = 1.417124
$ws.Cells.Item(2, 9).Value = 0.02676815992729067
$ws.Cells.Item(2, 10).Value = 0.02676815992729067
$ws.Cells.Item(2, 13).Value = 0.8689623333333333
$ws.Cells.Item(2, 14).Value = 2.606887
$ws.Cells.Item(2, 15).Value = 0.2842399713021852
$ws.Cells.Item(2, 16).Value = 0.2842399713021851
$ws.Cells.Item(2, 17).Value = 0.4104757925542222
$ws.Cells.Item(2, 18).Value = 3.694282132988
$ws.Cells.Item(2, 19).Value = 0.007608581009545403
$ws.Cells.Item(2, 20).Value = 0.007608581009545401

$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.4723746666666667
$ws.Cells.Item(3, 8).Value = 1.417124
$ws.Cells.Item(3, 9).Value = 0.02676815992729067
$ws.Cells.Item(3, 10).Value = 0.02676815992729067
$ws.Cells.Item(3, 15).Value = 0.3353618792271216
$ws.Cells.Item(3, 16).Value = 0.3353618792271216
$ws.Cells.Item(3, 17).Value = 0.4843018120835555
$ws.Cells.Item(3, 18).Value = 4.358716308752
$ws.Cells.Item(3, 19).Value = 0.00897702041666833
$ws.Cells.Item(3, 20).Value = 0.008977020416668329

$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.4723746666666667
$ws.Cells.Item(4, 8).Value = 1.417124
$ws.Cells.Item(4, 9).Value = 0.02676815992729067
$ws.Cells.Item(4, 10).Value = 0.02676815992729067
$ws.Cells.Item(4, 13).Value = 0.378697
$ws.Cells.Item(4, 14).Value = 1.136091
$ws.Cells.Item(4, 15).Value = 0.1238728311724562
$ws.Cells.Item(4, 16).Value = 0.1238728311724562
$ws.Cells.Item(4, 17).Value = 0.1788868691426667
$ws.Cells.Item(4, 18).Value = 1.609981822284
$ws.Cells.Item(4, 19).Value = 0.003315847755470585
$ws.Cells.Item(4, 20).Value = 0.003315847755470584

$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.4723746666666667
$ws.Cells.Item(5, 8).Value = 1.417124
$ws.Cells.Item(5, 9).Value = 0.02676815992729067
$ws.Cells.Item(5, 10).Value = 0.02676815992729067
$ws.Cells.Item(5, 13).Value = 0.5162433333333333
$ws.Cells.Item(5, 14).Value = 1.54873
$ws.Cells.Item(5, 15).Value = 0.1688646154416487
$ws.Cells.Item(5, 16).Value = 0.1688646154416487
$ws.Cells.Item(5, 17).Value = 0.2438602725022222
$ws.Cells.Item(5, 18).Value = 2.19474245252
$ws.Cells.Item(5, 19).Value = 0.004520195032202489
$ws.Cells.Item(5, 20).Value = 0.004520195032202488

$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.4723746666666667
$ws.Cells.Item(6, 8).Value = 1.417124
$ws.Cells.Item(6, 9).Value = 0.02676815992729067
$ws.Cells.Item(6, 10).Value = 0.02676815992729067
$ws.Cells.Item(6, 13).Value = 0.2679913333333333
$ws.Cells.Item(6, 14).Value = 0.803974
$ws.Cells.Item(6, 15).Value = 0.08766070285658835
$ws.Cells.Item(6, 16).Value = 0.08766070285658834
$ws.Cells.Item(6, 17).Value = 0.1265923167528889
$ws.Cells.Item(6, 18).Value = 1.139330850776
$ws.Cells.Item(6, 19).Value = 0.002346515713403863
$ws.Cells.Item(6, 20).Value = 0.002346515713403862

$ws.Cells.Item(7, 9).Value = 0.9656838605972748
$ws.Cells.Item(7, 10).Value = 0.9656838605972748
$ws.Cells.Item(7, 13).Value = 0.8689623333333333
$ws.Cells.Item(7, 14).Value = 2.606887
$ws.Cells.Item(7, 15).Value = 0.2842399713021852
$ws.Cells.Item(7, 16).Value = 0.2842399713021851
$ws.Cells.Item(7, 17).Value = 14.80825910754366
$ws.Cells.Item(7, 18).Value = 133.274331967893
$ws.Cells.Item(7, 19).Value = 0.2744859528231528
$ws.Cells.Item(7, 20).Value = 0.2744859528231527

$ws.Cells.Item(8, 9).Value = 0.9656838605972748
$ws.Cells.Item(8, 10).Value = 0.9656838605972748
$ws.Cells.Item(8, 15).Value = 0.3353618792271216
$ws.Cells.Item(8, 16).Value = 0.3353618792271216
$ws.Cells.Item(8, 19).Value = 0.3238535542292039
$ws.Cells.Item(8, 20).Value = 0.3238535542292038

$ws.Cells.Item(9, 9).Value = 0.9656838605972748
$ws.Cells.Item(9, 10).Value = 0.9656838605972748
$ws.Cells.Item(9, 13).Value = 0.378697
$ws.Cells.Item(9, 14).Value = 1.136091
$ws.Cells.Item(9, 15).Value = 0.1238728311724562
$ws.Cells.Item(9, 16).Value = 0.1238728311724562
$ws.Cells.Item(9, 17).Value = 6.453494109160999
$ws.Cells.Item(9, 18).Value = 58.08144698244899
$ws.Cells.Item(9, 19).Value = 0.119621993829732
$ws.Cells.Item(9, 20).Value = 0.1196219938297319

$ws.Cells.Item(10, 9).Value = 0.9656838605972748
$ws.Cells.Item(10, 10).Value = 0.9656838605972748
$ws.Cells.Item(10, 13).Value = 0.5162433333333333
$ws.Cells.Item(10, 14).Value = 1.54873
$ws.Cells.Item(10, 15).Value = 0.1688646154416487
$ws.Cells.Item(10, 16).Value = 0.1688646154416487
$ws.Cells.Item(10, 17).Value = 8.797464227496665
$ws.Cells.Item(10, 18).Value = 79.17717804747
$ws.Cells.Item(10, 19).Value = 0.1630698337579655
$ws.Cells.Item(10, 20).Value = 0.1630698337579655

$ws.Cells.Item(11, 9).Value = 0.9656838605972748
$ws.Cells.Item(11, 10).Value = 0.9656838605972748
$ws.Cells.Item(11, 13).Value = 0.2679913333333333
$ws.Cells.Item(11, 14).Value = 0.803974
$ws.Cells.Item(11, 15).Value = 0.08766070285658835
$ws.Cells.Item(11, 16).Value = 0.08766070285658834
$ws.Cells.Item(11, 17).Value = 4.566924192620665
$ws.Cells.Item(11, 18).Value = 41.102317733586
$ws.Cells.Item(11, 19).Value = 0.0846525259572208
$ws.Cells.Item(11, 20).Value = 0.08465252595722078

$ws.Cells.Item(12, 7).Value = 0.1331983333333333
$ws.Cells.Item(12, 8).Value = 0.399595
$ws.Cells.Item(12, 9).Value = 0.007547979475434553
$ws.Cells.Item(12, 10).Value = 0.007547979475434553
$ws.Cells.Item(12, 13).Value = 0.8689623333333333
$ws.Cells.Item(12, 14).Value = 2.606887
$ws.Cells.Item(12, 15).Value = 0.2842399713021852
$ws.Cells.Item(12, 16).Value = 0.2842399713021851
$ws.Cells.Item(12, 17).Value = 0.1157443345294444
$ws.Cells.Item(12, 18).Value = 1.041699010765
$ws.Cells.Item(12, 19).Value = 0.002145437469487
$ws.Cells.Item(12, 20).Value = 0.002145437469487

$ws.Cells.Item(13, 7).Value = 0.1331983333333333
$ws.Cells.Item(13, 8).Value = 0.399595
$ws.Cells.Item(13, 9).Value = 0.007547979475434553
$ws.Cells.Item(13, 10).Value = 0.007547979475434553
$ws.Cells.Item(13, 15).Value = 0.3353618792271216
$ws.Cells.Item(13, 16).Value = 0.3353618792271216
$ws.Cells.Item(13, 17).Value = 0.1365615024511111
$ws.Cells.Item(13, 18).Value = 1.22905352206
$ws.Cells.Item(13, 19).Value = 0.002531304581249475
$ws.Cells.Item(13, 20).Value = 0.002531304581249475

$ws.Cells.Item(14, 7).Value = 0.1331983333333333
$ws.Cells.Item(14, 8).Value = 0.399595
$ws.Cells.Item(14, 9).Value = 0.007547979475434553
$ws.Cells.Item(14, 10).Value = 0.007547979475434553
$ws.Cells.Item(14, 13).Value = 0.378697
$ws.Cells.Item(14, 14).Value = 1.136091
$ws.Cells.Item(14, 15).Value = 0.1238728311724562
$ws.Cells.Item(14, 16).Value = 0.1238728311724562
$ws.Cells.Item(14, 17).Value = 0.05044180923833334
$ws.Cells.Item(14, 18).Value = 0.453976283145
$ws.Cells.Item(14, 19).Value = 0.000934989587253669
$ws.Cells.Item(14, 20).Value = 0.0009349895872536688

$ws.Cells.Item(15, 7).Value = 0.1331983333333333
$ws.Cells.Item(15, 8).Value = 0.399595
$ws.Cells.Item(15, 9).Value = 0.007547979475434553
$ws.Cells.Item(15, 10).Value = 0.007547979475434553
$ws.Cells.Item(15, 13).Value = 0.5162433333333333
$ws.Cells.Item(15, 14).Value = 1.54873
$ws.Cells.Item(15, 15).Value = 0.1688646154416487
$ws.Cells.Item(15, 16).Value = 0.1688646154416487
$ws.Cells.Item(15, 17).Value = 0.06876275159444444
$ws.Cells.Item(15, 18).Value = 0.61886476435
$ws.Cells.Item(15, 19).Value = 0.001274586651480713
$ws.Cells.Item(15, 20).Value = 0.001274586651480713

$ws.Cells.Item(16, 7).Value = 0.1331983333333333
$ws.Cells.Item(16, 8).Value = 0.399595
$ws.Cells.Item(16, 9).Value = 0.007547979475434553
$ws.Cells.Item(16, 10).Value = 0.007547979475434553
$ws.Cells.Item(16, 13).Value = 0.2679913333333333
$ws.Cells.Item(16, 14).Value = 0.803974
$ws.Cells.Item(16, 15).Value = 0.08766070285658835
$ws.Cells.Item(16, 16).Value = 0.08766070285658834
$ws.Cells.Item(16, 17).Value = 0.03569599894777777
$ws.Cells.Item(16, 18).Value = 0.32126399053
$ws.Cells.Item(16, 19).Value = 0.000661661185963696
$ws.Cells.Item(16, 20).Value = 0.0006616611859636959
